# edit.ps1 - reproduces the OOXML diff via PowerPoint COM automation.
#
# The commit applies two independent changes:
#   1. The single table on the deck gets a new table style GUID
#      ({A9AA8743-7688-4006-9893-9C7A0D41B3A1} -> {5C3712F3-D4CE-4557-8D77-4246EB51600F}).
#   2. The presentation's theme ("Integral") has its 12 theme colors replaced
#      with the classic "Office Theme" palette (name/clrScheme-name attributes
#      are not exposed for writing by the PowerPoint object model, so only the
#      actual RGB values - the part that IS scriptable - are changed).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Re-style the table (search every slide/shape so this is resilient to
#    shape ordering instead of hard-coding slide/shape indices).
# ---------------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{5C3712F3-D4CE-4557-8D77-4246EB51600F}")
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Swap the theme palette from "Integral" to the stock "Office Theme"
#    colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in that order).
# ---------------------------------------------------------------------------
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeRGB.Length; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
